# Bug Report.docx update:
#  - resize the table's grid columns
#  - clear out stale handler/status text on bug 5
#  - append two new bug rows (6 and 7)
#  - move the "_GoBack" bookmark from the old last row onto the new last row

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Resize the table grid columns (dxa -> points, 20 twips per point) ---
$t.Columns(1).Width = 701 / 20
$t.Columns(2).Width = 2900 / 20
$t.Columns(3).Width = 2208 / 20
$t.Columns(4).Width = 2305 / 20
$t.Columns(5).Width = 1236 / 20

# --- 2. Drop the stray "_GoBack" bookmark sitting on bug 3's status cell ---
$d.Bookmarks("_GoBack").Delete()

# --- 3. Bug 5 ("Stack failure...") gains a handler + status ---
$t.Cell(6, 4).Range.Text = "Roland"
$t.Cell(6, 5).Range.Text = "Fixing."

# --- 4. Append bug 6 ---
$row6 = $t.Rows.Add()
$t.Cell($row6.Index, 1).Range.Text = "6."
$t.Cell($row6.Index, 2).Range.Text = "Masks bought during gameplay are not reflected in endofdayState"
$t.Cell($row6.Index, 3).Range.Text = "Yes"
$t.Cell($row6.Index, 5).Range.Text = "Unresolved"

# --- 5. Append bug 7 ---
$row7 = $t.Rows.Add()
$t.Cell($row7.Index, 1).Range.Text = "7."
$t.Cell($row7.Index, 2).Range.Text = "Cannot move camera when paused during game"
$t.Cell($row7.Index, 3).Range.Text = "Yes"
$t.Cell($row7.Index, 5).Range.Text = "Unresolved"

# --- 6. Re-plant "_GoBack" at the very end of the table (end of bug 7's status cell) ---
$lastCell = $t.Cell($row7.Index, 5)
$endRange = $lastCell.Range
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRange)
